$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Put the path to the images under the corresponding column."
$ws.Range("A2").Value = "The path must be the remaining path after the root folder, which is:"
$ws.Range("A3").Value = "/home/masoud/Documents/four-polar/fourPolar-io/target/test-classes/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel"
$ws.Range("A4").Value = "The files in each row must correspond to different polarizations of same sample."
